# Weekly data refresh: a new pair of rows (Primera/Segunda quality) for a
# newer price-reporting date is inserted at the top of the data block
# (row 82), pushing all the existing weekly observations down by two rows.
# The new rows reuse the same Volumen/Precio/Unidad/Origen values that used
# to sit in the old row 82/83, only the Fecha (date) changes to the new
# reporting date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 82; everything below (old rows 82-201) moves
# down to rows 84-203.
$ws.Rows("82:83").Insert()

# The old row 82/83 payload now lives at rows 84/85 - duplicate it back into
# the freshly inserted rows 82/83 so every column (Mercado ID, Mercado,
# Región, Codreg, Categoría ID/Categoría, Variedad, Calidad, Volumen,
# Precio mínimo/máximo/promedio, Unidad, Origen, Precio $/Kg, Kg o Unidades,
# Clasificación) keeps the same value as before.
$src = $ws.Range("A84:R85")
$dst = $ws.Range("A82:R83")
$src.Copy($dst)

# Only the report date actually changes for the newly-inserted pair.
$ws.Cells.Item(82, 4).Value = 44467
$ws.Cells.Item(83, 4).Value = 44467
